$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 548.5
$ws.Range("I11").Value = 548.5
$ws.Range("K11").Value = 548.5
$ws.Range("M11").Value = -408.5

$ws.Range("H12").Value = 105.22222
$ws.Range("J12").Value = 56
$ws.Range("L12").Value = 56
$ws.Range("N12").Value = -396

$ws.Range("H55").Value = 1726.9231
$ws.Range("J55").Value = 5216.75
$ws.Range("L55").Value = 5216.75
$ws.Range("N55").Value = -5644.75

$ws.Range("H69").Value = 33203
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 33203
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 99609
$ws.Range("N69").Value = -101357
$ws.Range("M69").ClearContents()

$ws.Range("H70").Value = 2362.111
$ws.Range("J70").Value = 2746.1667
$ws.Range("L70").Value = 8238.500100000001
$ws.Range("N70").Value = -8778.500100000001

$ws.Range("H72").Value = 33203
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 33203
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 298827
$ws.Range("N72").Value = -307563
$ws.Range("M72").ClearContents()

$ws.Range("H73").Value = 2362.111
$ws.Range("J73").Value = 2746.1667
$ws.Range("L73").Value = 8238.500100000001
$ws.Range("N73").Value = -10110.5001

$ws.Range("H100").Value = 3483.3845
$ws.Range("J100").Value = 6600.6
$ws.Range("L100").Value = 6600.6
$ws.Range("N100").Value = -7682.6

$ws.Range("H103").Value = 850.8333
$ws.Range("J103").Value = 901
$ws.Range("L103").Value = 2703
$ws.Range("N103").Value = -3875

$ws.Range("H133").Value = 73999.5
$ws.Range("J133").Value = 73999.5
$ws.Range("L133").Value = 73999.5
$ws.Range("N133").Value = -84119.5

$ws.Range("H138").Value = 2126.5588
$ws.Range("I138").Value = 812.17645
$ws.Range("J138").Value = 3440.9412
$ws.Range("K138").Value = 2436.52935
$ws.Range("L138").Value = 10322.8236
$ws.Range("M138").Value = 2703.47065
$ws.Range("N138").Value = -20602.8236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1268.2142
$ws.Range("I2").Value = 1346.75
$ws.Range("J2").Value = 797
$ws.Range("K2").Value = 1346.75
$ws.Range("L2").Value = 797
$ws.Range("M2").Value = -1233.75
$ws.Range("N2").Value = -1023

$ws.Range("H32").Value = 22739116
$ws.Range("I32").Value = 23821846
$ws.Range("K32").Value = 23821846
$ws.Range("M32").Value = -23821559

$ws.Range("H101").Value = 89969.71000000001
$ws.Range("J101").Value = 89969.71000000001
$ws.Range("L101").Value = 89969.71000000001
$ws.Range("N101").Value = -96459.71000000001

$ws.Range("H116").Value = 1268.2142
$ws.Range("I116").Value = 1346.75
$ws.Range("J116").Value = 797
$ws.Range("K116").Value = 1346.75
$ws.Range("L116").Value = 797
$ws.Range("M116").Value = 947.25
$ws.Range("N116").Value = -5385

$ws.Range("H130").Value = 105998.5
$ws.Range("J130").Value = 105998.5
$ws.Range("L130").Value = 105998.5
$ws.Range("N130").Value = -116038.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1268.2142
$ws.Range("I3").Value = 1346.75
$ws.Range("J3").Value = 797
$ws.Range("K3").Value = 1346.75
$ws.Range("L3").Value = 797
$ws.Range("M3").Value = -1232.75
$ws.Range("N3").Value = -1025

$ws.Range("H134").Value = 71318.07000000001
$ws.Range("I134").Value = 1395.7
$ws.Range("J134").Value = 211162.8
$ws.Range("K134").Value = 4187.1
$ws.Range("L134").Value = 633488.3999999999
$ws.Range("M134").Value = -1652.1
$ws.Range("N134").Value = -638558.3999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 200
$ws.Range("I12").Value = 200
$ws.Range("K12").Value = 200
$ws.Range("M12").Value = -30

$ws.Range("H19").Value = 500
$ws.Range("I19").Value = 500
$ws.Range("K19").Value = 500
$ws.Range("M19").Value = -330

$ws.Range("H24").Value = 500
$ws.Range("I24").Value = 500
$ws.Range("K24").Value = 500
$ws.Range("M24").Value = -330

$ws.Range("H31").Value = 691639
$ws.Range("I31").Value = 8598.75
$ws.Range("J31").Value = 1667410.8
$ws.Range("K31").Value = 8598.75
$ws.Range("L31").Value = 1667410.8
$ws.Range("M31").Value = -8303.75
$ws.Range("N31").Value = -1668000.8

$ws.Range("H34").Value = 691639
$ws.Range("I34").Value = 8598.75
$ws.Range("J34").Value = 1667410.8
$ws.Range("K34").Value = 8598.75
$ws.Range("L34").Value = 1667410.8
$ws.Range("M34").Value = -8396.75
$ws.Range("N34").Value = -1667814.8

$ws.Range("H43").Value = 47664
$ws.Range("J43").Value = 47664
$ws.Range("L43").Value = 47664
$ws.Range("N43").Value = -48032

$ws.Range("H101").Value = 47664
$ws.Range("J101").Value = 47664
$ws.Range("L101").Value = 47664
$ws.Range("N101").Value = -54154

$ws.Range("H107").Value = 850.34784
$ws.Range("I107").Value = 587.6111
$ws.Range("K107").Value = 587.6111
$ws.Range("M107").Value = 1332.3889

$ws.Range("H108").Value = 75339.336
$ws.Range("J108").Value = 75339.336
$ws.Range("L108").Value = 75339.336
$ws.Range("N108").Value = -83019.336

$ws.Range("H115").Value = 40976.332
$ws.Range("J115").Value = 40976.332
$ws.Range("L115").Value = 40976.332
$ws.Range("N115").Value = -43326.332

$ws.Range("H120").Value = 25666.666
$ws.Range("J120").Value = 25666.666
$ws.Range("L120").Value = 25666.666
$ws.Range("N120").Value = -32924.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 9.5
$ws.Range("I10").Value = 9.444445
$ws.Range("K10").Value = 28.333335
$ws.Range("M10").Value = 110.666665

$ws.Range("H11").Value = 621
$ws.Range("I11").Value = 505.2
$ws.Range("K11").Value = 1515.6
$ws.Range("M11").Value = -1375.6

$ws.Range("H46").Value = 687.625
$ws.Range("J46").Value = 1043
$ws.Range("L46").Value = 3129
$ws.Range("N46").Value = -3311

$ws.Range("H104").Value = 4748.75
$ws.Range("I104").Value = 4497.5
$ws.Range("K104").Value = 13492.5
$ws.Range("M104").Value = -10871.5

$ws.Range("H113").Value = 1361.5
$ws.Range("I113").Value = 659.6667
$ws.Range("J113").Value = 1782.6
$ws.Range("K113").Value = 1979.0001
$ws.Range("L113").Value = 5347.799999999999
$ws.Range("M113").Value = 190.9999
$ws.Range("N113").Value = -9687.799999999999

$ws.Range("H115").Value = 37153
$ws.Range("J115").Value = 70031
$ws.Range("L115").Value = 210093
$ws.Range("N115").Value = -212443

$ws.Range("H122").Value = 1802.3334
$ws.Range("J122").Value = 2463.4
$ws.Range("L122").Value = 22170.6
$ws.Range("N122").Value = -27070.6

$ws.Range("H124").Value = 1799.75
$ws.Range("I124").Value = 1799.75
$ws.Range("K124").Value = 5399.25
$ws.Range("M124").Value = -489.25

$ws.Range("H138").Value = 2041.3334
$ws.Range("I138").Value = 2237.25
$ws.Range("J138").Value = 1649.5
$ws.Range("K138").Value = 6711.75
$ws.Range("L138").Value = 4948.5
$ws.Range("M138").Value = -1571.75
$ws.Range("N138").Value = -15228.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 25006
$ws.Range("I23").Value = 25006
$ws.Range("K23").Value = 25006
$ws.Range("M23").Value = -24776

$ws.Range("H40").Value = 3981.2222
$ws.Range("I40").Value = 2695.5715
$ws.Range("J40").Value = 4799.364
$ws.Range("K40").Value = 2695.5715
$ws.Range("L40").Value = 4799.364
$ws.Range("M40").Value = -2559.5715
$ws.Range("N40").Value = -5071.364

$ws.Range("H43").Value = 2538247.5
$ws.Range("J43").Value = 35989.5
$ws.Range("L43").Value = 35989.5
$ws.Range("N43").Value = -36375.5

$ws.Range("H46").Value = 3007.9167
$ws.Range("J46").Value = 4521.25
$ws.Range("L46").Value = 4521.25
$ws.Range("N46").Value = -4897.25

$ws.Range("H100").Value = 3056.7222
$ws.Range("I100").Value = 3287.2144
$ws.Range("J100").Value = 2250
$ws.Range("K100").Value = 3287.2144
$ws.Range("L100").Value = 2250
$ws.Range("M100").Value = -2746.2144
$ws.Range("N100").Value = -3332

Write-Output "Applied all Behemoth_Profits price updates"
